$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$pairs = @(
    @("38+9=", "81-54="),
    @("70-46=", "96-27="),
    @("37+49=", "71-35="),
    @("56+28=", "65+17="),
    @("41-12=", "18+8="),
    @("65-58=", "58+23="),
    @("26+39=", "7+4="),
    @("40-35=", "92-33="),
    @("84-27=", "44-18="),
    @("43-25=", "81-75="),
    @("15+29=", "62-48="),
    @("5+7=", "30-13="),
    @("58-19=", "47-28="),
    @("55-6=", "94-25="),
    @("64+18=", "29+28="),
    @("36+35=", "58+15="),
    @("74-46=", "29+69="),
    @("7+76=", "36-19="),
    @("46-38=", "37-28="),
    @("61-28=", "9+53="),
    @("98-69=", "92-25="),
    @("54-9=", "46+5="),
    @("28+64=", "62-8="),
    @("83-79=", "84-48="),
    @("83-75=", "8+47="),
    @("54-16=", "79+7="),
    @("63+18=", "54-7="),
    @("91-19=", "17+24="),
    @("22-16=", "39+45="),
    @("25-18=", "40-27="),
    @("90-12=", "41-8="),
    @("94-8=", "19+23="),
    @("35+49=", "29+48="),
    @("93-28=", "94-48="),
    @("69+22=", "84-76="),
    @("59+18=", "89+3="),
    @("66+16=", "53+8="),
    @("55+36=", "34+38="),
    @("57+38=", "78-9="),
    @("43-24=", "4+27="),
    @("62-15=", "8+5="),
    @("92-23=", "27+16="),
    @("9+59=", "73-6="),
    @("7+85=", "92-16="),
    @("25+49=", "80-53="),
    @("64-16=", "50-25="),
    @("50-4=", "52-48="),
    @("42+39=", "50-45="),
    @("84-67=", "81-17="),
    @("40-23=", "71-45="),
    @("50-8=", "86-7="),
    @("94-57=", "29+69="),
    @("39+13=", "90-67="),
    @("45+48=", "16+35="),
    @("45+9=", "81-8="),
    @("9+39=", "65-56="),
    @("42-9=", "22-15="),
    @("80-19=", "93-45="),
    @("27+27=", "61-35="),
    @("71-24=", "17+14="),
    @("28+67=", "8+64="),
    @("46+25=", "53-16="),
    @("6+9=", "25-19="),
    @("91-26=", "31-23="),
    @("9+75=", "9+35="),
    @("98-19=", "65-37="),
    @("33+28=", "92-64="),
    @("69+24=", "31-22="),
    @("93-17=", "91-24="),
    @("14+17=", "79+2="),
    @("94-89=", "76+7="),
    @("64-35=", "12+9="),
    @("6+76=", "83-5="),
    @("38+19=", "8+55="),
    @("54-26=", "70-15="),
    @("52-9=", "84-17="),
    @("35-9=", "64+9="),
    @("64-17=", "50-31="),
    @("9+16=", "53-29="),
    @("91-19=", "61-36="),
    @("37-19=", "90-63="),
    @("28+56=", "7+57="),
    @("28+16=", "65-17="),
    @("15+58=", "9+66="),
    @("66-8=", "84-55="),
    @("74-25=", "27+69="),
    @("65-8=", "44+18="),
    @("62-24=", "45+36="),
    @("9+5=", "15+79="),
    @("48+46=", "91-25="),
    @("34+8=", "36-19="),
    @("78+13=", "4+49="),
    @("95-26=", "36+48="),
    @("25+29=", "44+28="),
    @("27+37=", "74+8="),
    @("19+4=", "38+28="),
    @("60-27=", "75+9="),
    @("78+3=", "63-14="),
    @("27+54=", "85-18="),
    @("88-19=", "60-55=")
)

for ($i = 0; $i -lt $pairs.Count; $i++) {
    $row = [int]([math]::Floor($i / 5)) + 1
    $col = ($i % 5) + 1
    $old = $pairs[$i][0]
    $new = $pairs[$i][1]
    $cell = $t.Cell($row, $col)
    $cell.Range.Text = $new
}
